$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (122) down to the
# three new rows (123-125) so the new cells pick up the same cell styles
# (bold/border index style for column A, date-number-format style for column E).
$ws.Range("A122:V122").Copy() | Out-Null
$ws.Range("A123:V125").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(123,1).Value = 122
$ws.Cells.Item(123,2).Value = 'costa-rica'
$ws.Cells.Item(123,3).Value = 'primera-division'
$ws.Cells.Item(123,4).Value = '2023-2024'
$ws.Cells.Item(123,5).Value = 45252.91666666666
$ws.Cells.Item(123,6).Value = 'Guanacasteca'
$ws.Cells.Item(123,7).Value = 4
$ws.Cells.Item(123,8).Value = 'Liberia'
$ws.Cells.Item(123,9).Value = 2
$ws.Cells.Item(123,10).Value = 1.81
$ws.Cells.Item(123,11).Value = '15/11/2023 22:12'
$ws.Cells.Item(123,12).Value = 2.05
$ws.Cells.Item(123,13).Value = '22/11/2023 21:53'
$ws.Cells.Item(123,14).Value = 3.66
$ws.Cells.Item(123,15).Value = '15/11/2023 22:12'
$ws.Cells.Item(123,16).Value = 3.59
$ws.Cells.Item(123,17).Value = '22/11/2023 21:57'
$ws.Cells.Item(123,18).Value = 4.28
$ws.Cells.Item(123,19).Value = '15/11/2023 22:12'
$ws.Cells.Item(123,20).Value = 3.57
$ws.Cells.Item(123,21).Value = '22/11/2023 21:57'
$ws.Cells.Item(123,22).Value = 'https://www.betexplorer.com/football/costa-rica/primera-division/guanacasteca-liberia/OO4bBtnM/'

$ws.Cells.Item(124,1).Value = 123
$ws.Cells.Item(124,2).Value = 'costa-rica'
$ws.Cells.Item(124,3).Value = 'primera-division'
$ws.Cells.Item(124,4).Value = '2023-2024'
$ws.Cells.Item(124,5).Value = 45253.08333333334
$ws.Cells.Item(124,6).Value = 'Puntarenas FC'
$ws.Cells.Item(124,7).Value = 1
$ws.Cells.Item(124,8).Value = 'Alajuelense'
$ws.Cells.Item(124,9).Value = 2
$ws.Cells.Item(124,10).Value = 4.04
$ws.Cells.Item(124,11).Value = '15/11/2023 18:12'
$ws.Cells.Item(124,12).Value = 4.2
$ws.Cells.Item(124,13).Value = '23/11/2023 01:59'
$ws.Cells.Item(124,14).Value = 3.58
$ws.Cells.Item(124,15).Value = '15/11/2023 18:12'
$ws.Cells.Item(124,16).Value = 3.74
$ws.Cells.Item(124,17).Value = '23/11/2023 01:59'
$ws.Cells.Item(124,18).Value = 1.88
$ws.Cells.Item(124,19).Value = '15/11/2023 18:12'
$ws.Cells.Item(124,20).Value = 1.85
$ws.Cells.Item(124,21).Value = '23/11/2023 01:59'
$ws.Cells.Item(124,22).Value = 'https://www.betexplorer.com/football/costa-rica/primera-division/puntarenas-fc-alajuelense/YcfY4vPk/'

$ws.Cells.Item(125,1).Value = 124
$ws.Cells.Item(125,2).Value = 'costa-rica'
$ws.Cells.Item(125,3).Value = 'primera-division'
$ws.Cells.Item(125,4).Value = '2023-2024'
$ws.Cells.Item(125,5).Value = 45253.125
$ws.Cells.Item(125,6).Value = 'Herediano'
$ws.Cells.Item(125,7).Value = 3
$ws.Cells.Item(125,8).Value = 'AD Santos'
$ws.Cells.Item(125,9).Value = 0
$ws.Cells.Item(125,10).Value = 1.32
$ws.Cells.Item(125,11).Value = '16/11/2023 05:12'
$ws.Cells.Item(125,12).Value = 1.36
$ws.Cells.Item(125,13).Value = '23/11/2023 02:57'
$ws.Cells.Item(125,14).Value = 5.14
$ws.Cells.Item(125,15).Value = '16/11/2023 05:12'
$ws.Cells.Item(125,16).Value = 5.03
$ws.Cells.Item(125,17).Value = '23/11/2023 02:58'
$ws.Cells.Item(125,18).Value = 9.08
$ws.Cells.Item(125,19).Value = '16/11/2023 05:12'
$ws.Cells.Item(125,20).Value = 8.73
$ws.Cells.Item(125,21).Value = '23/11/2023 02:58'
$ws.Cells.Item(125,22).Value = 'https://www.betexplorer.com/football/costa-rica/primera-division/herediano-santos-de-guapiles/xS82A01S/'

Write-Output "Inserted rows 123-125"
